$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Format the price (D-column) cells as text first so numeric-looking strings
# are preserved exactly (with trailing zeros etc.) instead of being converted
# to floating point numbers by Excel. (Applied per-cell because NumberFormat
# on a multi-area Range object only affects the first area.)
$dRefs = @("D2", "D4", "D5", "D6", "D7", "D8", "D9", "D10", "D11", "D12", "D13", "D14", "D15", "D16", "D18", "D20", "D21", "D22", "D23", "D25", "D26", "D27", "D28", "D40", "D41", "D42", "D43", "D44", "D45", "D48")
foreach ($ref in $dRefs) {
    $ws.Range($ref).NumberFormat = "@"
}

# Update price values (Column D)
$ws.Range("D2").Value = "245.61"
$ws.Range("D4").Value = "5.130"
$ws.Range("D5").Value = "0.05590"
$ws.Range("D6").Value = "6.485"
$ws.Range("D7").Value = "3.028"
$ws.Range("D8").Value = "0.8187"
$ws.Range("D9").Value = "0.8481"
$ws.Range("D10").Value = "0.1341"
$ws.Range("D11").Value = "0.06953"
$ws.Range("D12").Value = "0.02869"
$ws.Range("D13").Value = "0.09393"
$ws.Range("D14").Value = "0.001517"
$ws.Range("D15").Value = "0.0005980"
$ws.Range("D16").Value = "0.006087"
$ws.Range("D18").Value = "2.118"
$ws.Range("D20").Value = "0.03234"
$ws.Range("D21").Value = "0.1320"
$ws.Range("D22").Value = "3.742"
$ws.Range("D23").Value = "0.04687"
$ws.Range("D25").Value = "0.001248"
$ws.Range("D26").Value = "0.004602"
$ws.Range("D27").Value = "0.00009599"
$ws.Range("D28").Value = "0.0001390"
$ws.Range("D40").Value = "0.03655"
$ws.Range("D41").Value = "0.1364"
$ws.Range("D42").Value = "0.006110"
$ws.Range("D43").Value = "0.002465"
$ws.Range("D44").Value = "0.007782"
$ws.Range("D45").Value = "0.00005305"
$ws.Range("D48").Value = "0.002127"

# Update coin name / link / label values (Columns B, C, E)
$ws.Range("B10").Value = "WazirX"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("B13").Value = "BitMartToken"
$ws.Range("B14").Value = "BitForexToken"
$ws.Range("B15").Value = "One"
$ws.Range("B41").Value = "BKEXToken"
$ws.Range("B42").Value = "KickToken"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("C15").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("C41").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("C42").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("E10").Value = "9WazirXWRX"
$ws.Range("E11").Value = "10MandalaExchangeTokenMDX"
$ws.Range("E12").Value = "11BitrueCoinBTR"
$ws.Range("E13").Value = "12BitMartTokenBMX"
$ws.Range("E14").Value = "13BitForexTokenBF"
$ws.Range("E15").Value = "14OneONE"
$ws.Range("E41").Value = "40BKEXTokenBKKBestin24h"
$ws.Range("E42").Value = "41KickTokenKICK"
